# Sprint Review Protocol: add review minutes ("Mitschrift von Review hinzugefuegt")
# Mark requirements 1-3 (rows 16-18) as "Passed" (column F) with "ja",
# and move the active selection to F19 (next empty row in that column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = "ja"
$ws.Range("F17").Value = "ja"
$ws.Range("F18").Value = "ja"

$ws.Range("F19").Select()
